# test2.xlsx: fix up the A1 text and leave the selection on A2, matching
# the author's interactive edit (typed the rest of the sentence into A1,
# then moved off the cell to A2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "test2 this is good"
$ws.Range("A2").Select()

